# =====================================================================
# Update workbook per commit "Actualizacion automatica 2025-11-10 16:30:09"
# =====================================================================
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$currencyFmt = "`"$`"#,##0.00"
$pctFmt = "0.00%"
$xlRight = -4152

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" - individual cell updates
# ---------------------------------------------------------------------
$ws1.Range("L66").Value2 = 89.56
$ws1.Range("D69").Value2 = 1852.03
$ws1.Range("M98").Value2 = 149.69
$ws1.Range("M129").Value2 = 56.86
$ws1.Range("L350").Value2 = "9 de 348"
$ws1.Range("M350").Value2 = "21 de 348"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" - individual cell updates
# ---------------------------------------------------------------------
$ws2.Range("F66").Value2 = 96.69
$ws2.Range("F69").Value2 = 1445.12
$ws2.Range("F98").Value2 = 149.69
$ws2.Range("F129").Value2 = 56.86
$ws2.Range("F354").Value2 = 24812.42

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" - rebuild with full detail data
# ---------------------------------------------------------------------
$sheet3Data = @(
    @{ Row=2; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="240X120 PORCELANATO"; C=129.6; D=0; E=129.6; F=0; IsTotal=$false },
    @{ Row=3; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="240X80 PORCELANATO"; C=3592.51; D=0; E=3592.51; F=0; IsTotal=$false },
    @{ Row=4; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="FREGADEROS DE COCINA"; C=207.39; D=0; E=207.39; F=0; IsTotal=$false },
    @{ Row=5; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=6; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="INODOROS"; C=660.6; D=23.4; E=637.2; F=0.03542234332425068; IsTotal=$false },
    @{ Row=7; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="LAVABOS"; C=93.90000000000001; D=0; E=93.90000000000001; F=0; IsTotal=$false },
    @{ Row=8; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="NO RESURTIBLES"; C=350; D=0; E=350; F=0; IsTotal=$false },
    @{ Row=9; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="OTROS"; C=0; D=0; E=0; F=0; IsTotal=$false },
    @{ Row=10; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="PANELES DECORATIVOS"; C=388.107983534392; D=0; E=388.107983534392; F=0; IsTotal=$false },
    @{ Row=11; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="PIEDRA SINTERIZADA"; C=3446; D=-142.56; E=3588.56; F=-0.04136970400464306; IsTotal=$false },
    @{ Row=12; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="PORCELANATO"; C=31214; D=1592.89; E=29621.11; F=0.05103126802075992; IsTotal=$false },
    @{ Row=13; HasA=$true; A="ALMEIDA CUATIN JHONATHANN CARLOS"; B="PUERTAS DE SEGURIDAD"; C=111.043665120341; D=0; E=111.043665120341; F=0; IsTotal=$false },
    @{ Row=14; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="240X120 PORCELANATO"; C=1837.54; D=0; E=1837.54; F=0; IsTotal=$false },
    @{ Row=15; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="240X80 PORCELANATO"; C=14825.41; D=3240.38; E=11585.03; F=0.2185693346760731; IsTotal=$false },
    @{ Row=16; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="FREGADEROS DE COCINA"; C=789.38; D=220.48; E=568.9; F=0.2793078111935949; IsTotal=$false },
    @{ Row=17; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=18; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="INODOROS"; C=1605; D=489.6; E=1115.4; F=0.3050467289719626; IsTotal=$false },
    @{ Row=19; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="LAVABOS"; C=886.711016287574; D=136; E=750.711016287574; F=0.1533757870398366; IsTotal=$false },
    @{ Row=20; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="NO RESURTIBLES"; C=350; D=0; E=350; F=0; IsTotal=$false },
    @{ Row=21; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="OTROS"; C=0; D=0; E=0; F=0; IsTotal=$false },
    @{ Row=22; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="PANELES DECORATIVOS"; C=9916; D=-406.91; E=10322.91; F=-0.04103569987898346; IsTotal=$false },
    @{ Row=23; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="PIEDRA SINTERIZADA"; C=16148; D=917.5; E=15230.5; F=0.05681818181818182; IsTotal=$false },
    @{ Row=24; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="PORCELANATO"; C=50307; D=251.59; E=50055.41; F=0.005001093287216491; IsTotal=$false },
    @{ Row=25; HasA=$true; A="CASTRO ALCIVAR EDA MARIA"; B="PUERTAS DE SEGURIDAD"; C=1110.43665120341; D=124.22; E=986.21665120341; F=0.1118659041606556; IsTotal=$false },
    @{ Row=26; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="240X120 PORCELANATO"; C=2826.66; D=0; E=2826.66; F=0; IsTotal=$false },
    @{ Row=27; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="240X80 PORCELANATO"; C=6623.26; D=1866.24; E=4757.02; F=0.2817706084314975; IsTotal=$false },
    @{ Row=28; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="FREGADEROS DE COCINA"; C=844.7; D=0; E=844.7; F=0; IsTotal=$false },
    @{ Row=29; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=30; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="INODOROS"; C=2907.58368146026; D=0; E=2907.58368146026; F=0; IsTotal=$false },
    @{ Row=31; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="LAVABOS"; C=1320; D=23.4; E=1296.6; F=0.01772727272727273; IsTotal=$false },
    @{ Row=32; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="NO RESURTIBLES"; C=415; D=0; E=415; F=0; IsTotal=$false },
    @{ Row=33; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="OTROS"; C=0; D=0; E=0; F=0; IsTotal=$false },
    @{ Row=34; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="PANELES DECORATIVOS"; C=4312; D=405.57; E=3906.43; F=0.0940561224489796; IsTotal=$false },
    @{ Row=35; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="PIEDRA SINTERIZADA"; C=14235.99; D=4962.1; E=9273.889999999999; F=0.3485602336051093; IsTotal=$false },
    @{ Row=36; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="PORCELANATO"; C=64944; D=2869.79; E=62074.21; F=0.04418868563685636; IsTotal=$false },
    @{ Row=37; HasA=$true; A="GUERRERO FAREZ FABIAN MAURICIO"; B="PUERTAS DE SEGURIDAD"; C=440.653177778119; D=0; E=440.653177778119; F=0; IsTotal=$false },
    @{ Row=38; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="240X120 PORCELANATO"; C=129.6; D=0; E=129.6; F=0; IsTotal=$false },
    @{ Row=39; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="240X80 PORCELANATO"; C=2564; D=-152.64; E=2716.64; F=-0.05953198127925116; IsTotal=$false },
    @{ Row=40; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="FREGADEROS DE COCINA"; C=207.39; D=0; E=207.39; F=0; IsTotal=$false },
    @{ Row=41; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=42; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="INODOROS"; C=2907.58368146026; D=0; E=2907.58368146026; F=0; IsTotal=$false },
    @{ Row=43; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="LAVABOS"; C=383.4; D=0; E=383.4; F=0; IsTotal=$false },
    @{ Row=44; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="NO RESURTIBLES"; C=415; D=0; E=415; F=0; IsTotal=$false },
    @{ Row=45; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="OTROS"; C=0; D=0; E=0; F=0; IsTotal=$false },
    @{ Row=46; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="PANELES DECORATIVOS"; C=1388; D=0; E=1388; F=0; IsTotal=$false },
    @{ Row=47; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="PIEDRA SINTERIZADA"; C=2678; D=-1151.4; E=3829.4; F=-0.4299477221807319; IsTotal=$false },
    @{ Row=48; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="PORCELANATO"; C=44418; D=3045.53; E=41372.47; F=0.06856522130667747; IsTotal=$false },
    @{ Row=49; HasA=$true; A="HIDALGO HIDALGO PEDRO GUSTAVO"; B="PUERTAS DE SEGURIDAD"; C=222.087330240682; D=0; E=222.087330240682; F=0; IsTotal=$false },
    @{ Row=50; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="240X120 PORCELANATO"; C=129.6; D=0; E=129.6; F=0; IsTotal=$false },
    @{ Row=51; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="240X80 PORCELANATO"; C=1867.69; D=0; E=1867.69; F=0; IsTotal=$false },
    @{ Row=52; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="FREGADEROS DE COCINA"; C=1987.7; D=2172.1; E=-184.3999999999999; F=1.092770538813704; IsTotal=$false },
    @{ Row=53; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=54; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="INODOROS"; C=1815; D=0; E=1815; F=0; IsTotal=$false },
    @{ Row=55; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="LAVABOS"; C=383.4; D=0; E=383.4; F=0; IsTotal=$false },
    @{ Row=56; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="NO RESURTIBLES"; C=415; D=29.49; E=385.51; F=0.07106024096385542; IsTotal=$false },
    @{ Row=57; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="OTROS"; C=0; D=0; E=0; F=0; IsTotal=$false },
    @{ Row=58; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="PANELES DECORATIVOS"; C=388.107983534392; D=0; E=388.107983534392; F=0; IsTotal=$false },
    @{ Row=59; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="PIEDRA SINTERIZADA"; C=1440.92; D=0; E=1440.92; F=0; IsTotal=$false },
    @{ Row=60; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="PORCELANATO"; C=48041; D=1602.63; E=46438.37; F=0.03335963031577194; IsTotal=$false },
    @{ Row=61; HasA=$true; A="LINDAO ZUÑIGA BRYAN JOSE"; B="PUERTAS DE SEGURIDAD"; C=1332.52398144409; D=0; E=1332.52398144409; F=0; IsTotal=$false },
    @{ Row=62; HasA=$true; A="OFICINA-CATAECSA"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=63; HasA=$true; A="OFICINA-CATAECSA"; B="OTROS"; C=0; D=3241.06; E=-3241.06; F=0; IsTotal=$false },
    @{ Row=64; HasA=$true; A="OFICINA-CATAECSA"; B="PORCELANATO"; C=26000; D=8238.030000000001; E=17761.97; F=0.3168473076923077; IsTotal=$false },
    @{ Row=65; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="240X120 PORCELANATO"; C=129.6; D=0; E=129.6; F=0; IsTotal=$false },
    @{ Row=66; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="240X80 PORCELANATO"; C=2344.03; D=0; E=2344.03; F=0; IsTotal=$false },
    @{ Row=67; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="FREGADEROS DE COCINA"; C=207.39; D=0; E=207.39; F=0; IsTotal=$false },
    @{ Row=68; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="GRIFERIAS"; C=86.41; D=0; E=86.41; F=0; IsTotal=$false },
    @{ Row=69; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="INODOROS"; C=855.91; D=0; E=855.91; F=0; IsTotal=$false },
    @{ Row=70; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="LAVABOS"; C=383; D=0; E=383; F=0; IsTotal=$false },
    @{ Row=71; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="NO RESURTIBLES"; C=415; D=0; E=415; F=0; IsTotal=$false },
    @{ Row=72; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="OTROS"; C=0; D=0; E=0; F=0; IsTotal=$false },
    @{ Row=73; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="PANELES DECORATIVOS"; C=388.107983534392; D=0; E=388.107983534392; F=0; IsTotal=$false },
    @{ Row=74; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="PIEDRA SINTERIZADA"; C=902.88; D=443.44; E=459.44; F=0.4911394648236754; IsTotal=$false },
    @{ Row=75; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="PORCELANATO"; C=34701; D=317.92; E=34383.08; F=0.009161695628368058; IsTotal=$false },
    @{ Row=76; HasA=$true; A="RIOS CARRION ANGEL BENIGNO"; B="PUERTAS DE SEGURIDAD"; C=364.412605947529; D=0; E=364.412605947529; F=0; IsTotal=$false },
    @{ Row=77; HasA=$false; A=$null; B="TOTAL"; C=417248.6797415454; D=34359.85000000001; E=382888.8297415455; F=0.08234861287346286; IsTotal=$true }
)
foreach ($item in $sheet3Data) {
    $r = $item.Row
    if ($item.HasA) {
        $ws3.Cells.Item($r, 1).Value2 = $item.A
    }
    $ws3.Cells.Item($r, 2).Value2 = $item.B
    if ($item.IsTotal) {
        $ws3.Cells.Item($r, 2).HorizontalAlignment = $xlRight
    }
    $ws3.Cells.Item($r, 3).Value2 = $item.C
    $ws3.Cells.Item($r, 3).NumberFormat = $currencyFmt
    $ws3.Cells.Item($r, 4).Value2 = $item.D
    $ws3.Cells.Item($r, 4).NumberFormat = $currencyFmt
    $ws3.Cells.Item($r, 5).Value2 = $item.E
    $ws3.Cells.Item($r, 5).NumberFormat = $currencyFmt
    $ws3.Cells.Item($r, 6).Value2 = $item.F
    $ws3.Cells.Item($r, 6).NumberFormat = $pctFmt
}

# Column widths for CUMPLIMIENTO MENSUAL sheet (target stored width: 34,22,22,13,24,26)
$ws3.Columns.Item(1).ColumnWidth = 33.166666666666664
$ws3.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws3.Columns.Item(3).ColumnWidth = 21.166666666666668
$ws3.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws3.Columns.Item(5).ColumnWidth = 23.166666666666668
$ws3.Columns.Item(6).ColumnWidth = 25.166666666666668

Write-Host "Edit complete"
